$d = $word.ActiveDocument

# The existing "_GoBack" bookmark sits at the end of the first (only)
# paragraph; it will be re-created at the end of the new last paragraph
# below, so drop it here first.
$d.Bookmarks.Item("_GoBack").Delete()

# Append two new paragraphs after the existing "Hello world!" paragraph:
#   - "Hi" (eastAsia-hinted run) + " china" (plain run)
#   - "你好" (eastAsia-hinted run) followed by the relocated _GoBack bookmark
# Building this via raw WordOpenXML lets us control each run's formatting
# precisely (so "Hi"/" china" stay as two distinct runs, and the bookmark
# lands exactly at the end of the "你好" paragraph) instead of relying on
# the editor's ambient formatting-inheritance/bookmark-insertion quirks.
$r = $d.Range($d.Content.End, $d.Content.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Hi</w:t></w:r><w:r><w:t xml:space="preserve"> china</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>你好</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
